$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (post-edit) for rows 4-8, columns D, J, K, L, M, N, P, Q
# This corresponds to a cyclic re-ordering of the weekly price rows:
# new row4 <- old row6, new row5 <- old row7, new row6 <- old row8,
# new row7 <- old row4, new row8 <- old row5

$data = @{
    4 = @{ D = 44313; J = 20; K = 15000; L = 15000; M = 15000; N = "`$/caja 15 kilos empedrada"; P = 1000; Q = 15 }
    5 = @{ D = 44313; J = 20; K = 30000; L = 30000; M = 30000; N = "`$/caja 20 kilos empedrada"; P = 1500; Q = 20 }
    6 = @{ D = 44293; J = 10; K = 25000; L = 25000; M = 25000; N = "`$/caja 15 kilos empedrada"; P = 1667; Q = 15 }
    7 = @{ D = 44280; J = 30; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
    8 = @{ D = 44285; J = 20; K = 25000; L = 25000; M = 25000; N = "`$/caja 18 kilos empedrada"; P = 1389; Q = 18 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
}
